$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColToIndex([string]$col) {
    $result = 0
    for ($i = 0; $i -lt $col.Length; $i++) {
        $result = $result * 26 + ([int][char]$col[$i] - [int][char]'A' + 1)
    }
    return $result
}

$row2 = @{
    "D" = 1416
    "E" = 210
    "F" = 210
    "G" = 194
    "H" = 144
    "I" = 74
    "J" = 70
    "K" = 2604
    "L" = 360
    "M" = 2244
    "N" = 1912
    "O" = 332
    "P" = 58
    "Q" = 147
    "R" = -120
    "S" = 40
    "T" = 19
    "U" = 128
    "V" = 57
    "W" = 14.83
    "X" = 10.15
    "Y" = 3.99
    "Z" = 5.81
    "AA" = 16.03
    "AB" = 3125.71
    "AC" = 6372
    "AD" = 13.1
    "AE" = 166683
    "AF" = 0.5
    "AG" = 0
    "AH" = 0
    "AI" = 0
    "AJ" = 1154482
}

$row3 = @{
    "D" = 1390
    "E" = 227
    "F" = 227
    "G" = 361
    "H" = 296
    "I" = 216
    "J" = 80
    "K" = 2849
    "L" = 391
    "M" = 2459
    "N" = 2048
    "O" = 411
    "P" = 58
    "Q" = 268
    "R" = -116
    "S" = 29
    "T" = 26
    "U" = 243
    "V" = 88
    "W" = 16.33
    "X" = 21.3
    "Y" = 10.91
    "Z" = 10.86
    "AA" = 15.89
    "AB" = 3493.91
    "AC" = 18717
    "AD" = 4.24
    "AE" = 178497
    "AF" = 0.44
    "AG" = 0
    "AH" = 0
    "AI" = 0
    "AJ" = 1154482
}

$row4 = @{
    "D" = 1231
    "E" = 198
    "F" = 198
    "G" = 209
    "H" = 171
    "I" = 92
    "J" = 79
    "K" = 2939
    "L" = 320
    "M" = 2619
    "N" = 2139
    "O" = 480
    "P" = 58
    "Q" = 166
    "R" = -180
    "S" = -33
    "T" = 23
    "U" = 143
    "V" = 72
    "W" = 16.11
    "X" = 13.92
    "Y" = 4.41
    "Z" = 5.92
    "AA" = 12.23
    "AB" = 3659.16
    "AC" = 7995
    "AD" = 10.09
    "AE" = 186452
    "AF" = 0.43
    "AG" = 500
    "AH" = 0.62
    "AI" = 6.22
    "AJ" = 1154482
}

$row5 = @{
    "D" = 1210
    "E" = 170
    "F" = 170
    "G" = 173
    "H" = 143
    "I" = 59
    "J" = 84
    "K" = 3024
    "L" = 330
    "M" = 2694
    "N" = 2162
    "O" = 532
    "P" = 58
    "Q" = 169
    "R" = 4
    "S" = -42
    "T" = 26
    "U" = 143
    "V" = 82
    "W" = 14.08
    "X" = 11.79
    "Y" = 2.74
    "Z" = 4.78
    "AA" = 12.25
    "AB" = 3754.78
    "AC" = 5107
    "AD" = 15.71
    "AE" = 188452
    "AF" = 0.43
    "AG" = 500
    "AH" = 0.62
    "AI" = 9.73
    "AJ" = 1154482
}

$row6 = @{
    "D" = 1226
    "E" = 133
    "F" = 133
    "G" = 159
    "H" = 126
    "I" = 63
    "K" = 3044
    "L" = 297
    "M" = 2748
    "N" = 2187
    "P" = 58
    "Q" = 102
    "R" = -34
    "S" = -99
    "T" = 11
    "U" = 91
    "V" = 52
    "W" = 10.87
    "X" = 10.32
    "Y" = 2.89
    "Z" = 4.17
    "AA" = 10.8
    "AB" = 3851.47
    "AC" = 5435
    "AD" = 11.76
    "AE" = 190644
    "AF" = 0.34
    "AG" = 500
    "AH" = 0.78
    "AI" = 9.140000000000001
    "AJ" = 1154482
}

$rows = @{
    2 = $row2
    3 = $row3
    4 = $row4
    5 = $row5
    6 = $row6
}

foreach ($rowNum in $rows.Keys) {
    $rowData = $rows[$rowNum]
    foreach ($col in $rowData.Keys) {
        $colIdx = ColToIndex $col
        $ws.Cells.Item($rowNum, $colIdx).Value = $rowData[$col]
    }
}

# Rows 7-9: clear all data columns (D..AI), keeping only A, B, C index/label columns
$ws.Range("D7:AI9").ClearContents()

Write-Host "Edit applied successfully"
